$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.326.09"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.839.17"
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.18"
$ws.Range("E5").Value = "  -0.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6257"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07423"
$ws.Range("E8").Value = "  -0.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2894"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.77"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07719"
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.838.01"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.950"
$ws.Range("E13").Value = "  -1.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6733"
$ws.Range("E14").Value = "  -1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001020"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.73"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.219"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "29.381.84"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.73"
$ws.Range("E19").Value = "  +1.45%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.0000"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.337"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.09"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.465"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1345"
$ws.Range("E26").Value = "  -2.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.34"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07271"
$ws.Range("E28").Value = "  +13.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.454"
$ws.Range("E29").Value = "  +4.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.476"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.037"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.030"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("E33").Value = "  -0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.138"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6950"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.571"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01831"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.897"
$ws.Range("E38").Value = "  +4.24%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.815"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "1.229.60"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("E41").Value = "  +3.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9995"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").Value = "1.990.46"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.50"
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.42"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("E46").Value = "  +4.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.702"
$ws.Range("E47").Value = "  -2.99%  "
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1135"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.864"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3895"
$ws.Range("E51").Value = "  -1.42%  "
